# Signal table update:
#  - insert a new row at row 6 (ALUsrc / EX / source-of-Reg-output / source-of-immediate)
#  - shift the former rows 6-11 down to 7-12
#  - update window view + selection
#  - add a page setup (paper size / orientation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the book window a bit to the right (workbookView xWindow 1170 -> 3030)
$win = $excel.ActiveWindow
$win.Left = 3030
$win.Top = 105

# Insert a new row above the current row 6 ("ALUop" row), shifting the rest down.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new ALUsrc signal.
$ws.Range("A6").Value = "ALUsrc"
$ws.Range("B6").Value = "EX"
$ws.Range("C6").Value = "来自寄存器堆的输出"
$ws.Range("D6").Value = "来自符号扩展的立即数"

# Selection ends on D6, matching the authored selection in the workbook.
$ws.Range("D6").Select()

# Add page setup information (paper size + orientation) as in the target file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
